# Remove the "Rectangle 153" shape (SEDS USA Staff list) from slide 1.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $shape = $s.Shapes.Item($i)
    if ($shape.Name -eq "Rectangle 153") {
        $shape.Delete()
    }
}
